$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8107
$ws.Range("I62").Value = 6575.4
$ws.Range("J62").Value = 10659.667
$ws.Range("K62").Value = 6575.4
$ws.Range("L62").Value = 10659.667
$ws.Range("M62").Value = -5951.4
$ws.Range("N62").Value = -11907.667
$ws.Range("H65").Value = 8107
$ws.Range("I65").Value = 6575.4
$ws.Range("J65").Value = 10659.667
$ws.Range("K65").Value = 32877
$ws.Range("L65").Value = 53298.335
$ws.Range("M65").Value = -29757
$ws.Range("N65").Value = -59538.335
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -53736
$ws.Range("H82").Value = 3208.2
$ws.Range("I82").Value = 3208.2
$ws.Range("K82").Value = 9624.599999999999
$ws.Range("M82").Value = -9218.599999999999
$ws.Range("H85").Value = 3208.2
$ws.Range("I85").Value = 3208.2
$ws.Range("K85").Value = 9624.599999999999
$ws.Range("M85").Value = -8220.599999999999
$ws.Range("H100").Value = 2332
$ws.Range("I100").Value = 1358.2858
$ws.Range("J100").Value = 5740
$ws.Range("K100").Value = 1358.2858
$ws.Range("L100").Value = 5740
$ws.Range("M100").Value = -817.2858000000001
$ws.Range("N100").Value = -6822
$ws.Range("H111").Value = 1346.75
$ws.Range("I111").Value = 1346.75
$ws.Range("K111").Value = 4040.25
$ws.Range("M111").Value = -973.25
$ws.Range("H113").Value = 2923.5789
$ws.Range("I113").Value = 2915.1765
$ws.Range("J113").Value = 2995
$ws.Range("K113").Value = 2915.1765
$ws.Range("L113").Value = 2995
$ws.Range("M113").Value = 338.8235
$ws.Range("N113").Value = -9503
$ws.Range("H125").Value = 1350
$ws.Range("J125").Value = 1890
$ws.Range("L125").Value = 17010
$ws.Range("N125").Value = -21930
$ws.Range("H127").Value = 1907.7
$ws.Range("I127").Value = 1964.1111
$ws.Range("J127").Value = 1400
$ws.Range("K127").Value = 5892.3333
$ws.Range("L127").Value = 4200
$ws.Range("M127").Value = -932.3333000000002
$ws.Range("N127").Value = -14120
$ws.Range("H137").Value = 32227.303
$ws.Range("I137").Value = 751.3
$ws.Range("J137").Value = 45912.523
$ws.Range("K137").Value = 2253.9
$ws.Range("L137").Value = 137737.569
$ws.Range("M137").Value = 296.1000000000004
$ws.Range("N137").Value = -142837.569
$ws.Range("H138").Value = 3247.6
$ws.Range("I138").Value = 3434.5
$ws.Range("K138").Value = 10303.5
$ws.Range("M138").Value = -5163.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2115825.2
$ws.Range("I2").Value = 2585453
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 2585453
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -2585340
$ws.Range("N2").Value = -2726
$ws.Range("H32").Value = 2325.5635
$ws.Range("I32").Value = 1636.7858
$ws.Range("K32").Value = 1636.7858
$ws.Range("M32").Value = -1349.7858
$ws.Range("H61").Value = 3277.5173
$ws.Range("I61").Value = 2569.6
$ws.Range("J61").Value = 4036
$ws.Range("K61").Value = 2569.6
$ws.Range("L61").Value = 4036
$ws.Range("M61").Value = -2357.6
$ws.Range("N61").Value = -4460
$ws.Range("H86").Value = 24000
$ws.Range("I86").Value = 24000
$ws.Range("K86").Value = 24000
$ws.Range("M86").Value = -22814
$ws.Range("H89").Value = 24000
$ws.Range("I89").Value = 24000
$ws.Range("K89").Value = 72000
$ws.Range("M89").Value = -66072
$ws.Range("H116").Value = 2115825.2
$ws.Range("I116").Value = 2585453
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 2585453
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -2583159
$ws.Range("N116").Value = -7088
$ws.Range("H132").Value = 2457.1702
$ws.Range("I132").Value = 2332.6667
$ws.Range("J132").Value = 2864.6365
$ws.Range("K132").Value = 6998.000100000001
$ws.Range("L132").Value = 8593.9095
$ws.Range("M132").Value = -4468.000100000001
$ws.Range("N132").Value = -13653.9095
$ws.Range("H136").Value = 3277.5173
$ws.Range("I136").Value = 2569.6
$ws.Range("J136").Value = 4036
$ws.Range("K136").Value = 7708.799999999999
$ws.Range("L136").Value = 12108
$ws.Range("M136").Value = -5158.799999999999
$ws.Range("N136").Value = -17208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2115825.2
$ws.Range("I3").Value = 2585453
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 2585453
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -2585339
$ws.Range("N3").Value = -2728
$ws.Range("H20").Value = 2036.591
$ws.Range("I20").Value = 1863.5
$ws.Range("K20").Value = 1863.5
$ws.Range("M20").Value = -1616.5
$ws.Range("H134").Value = 3067.923
$ws.Range("I134").Value = 2820.6667
$ws.Range("J134").Value = 3624.25
$ws.Range("K134").Value = 8462.000100000001
$ws.Range("L134").Value = 10872.75
$ws.Range("M134").Value = -5927.000100000001
$ws.Range("N134").Value = -15942.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2899734
$ws.Range("I58").Value = 4348837
$ws.Range("J58").Value = 1528
$ws.Range("K58").Value = 4348837
$ws.Range("L58").Value = 1528
$ws.Range("M58").Value = -4348634
$ws.Range("N58").Value = -1934
$ws.Range("H99").Value = 2874.75
$ws.Range("I99").Value = 2374.5
$ws.Range("K99").Value = 2374.5
$ws.Range("M99").Value = -876.5
$ws.Range("H126").Value = 2874.75
$ws.Range("I126").Value = 2374.5
$ws.Range("K126").Value = 7123.5
$ws.Range("M126").Value = -4653.5
$ws.Range("H132").Value = 3115.7083
$ws.Range("J132").Value = 3677.4
$ws.Range("L132").Value = 11032.2
$ws.Range("N132").Value = -16092.2
$ws.Range("H136").Value = 2899734
$ws.Range("I136").Value = 4348837
$ws.Range("J136").Value = 1528
$ws.Range("K136").Value = 13046511
$ws.Range("L136").Value = 4584
$ws.Range("M136").Value = -13043961
$ws.Range("N136").Value = -9684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 249.13333
$ws.Range("J2").Value = 293.4
$ws.Range("L2").Value = 1760.4
$ws.Range("N2").Value = -1986.4
$ws.Range("H68").Value = 1827.46
$ws.Range("J68").Value = 1963.5227
$ws.Range("L68").Value = 5890.5681
$ws.Range("N68").Value = -7512.5681
$ws.Range("H71").Value = 1827.46
$ws.Range("J71").Value = 1963.5227
$ws.Range("L71").Value = 17671.7043
$ws.Range("N71").Value = -25783.7043
$ws.Range("H98").Value = 449.81818
$ws.Range("J98").Value = 443.75
$ws.Range("L98").Value = 1331.25
$ws.Range("N98").Value = -4327.25
$ws.Range("H107").Value = 1366.9231
$ws.Range("J107").Value = 1430.5714
$ws.Range("L107").Value = 4291.7142
$ws.Range("N107").Value = -8131.7142
$ws.Range("H131").Value = 11922625
$ws.Range("J131").Value = 18730.85
$ws.Range("L131").Value = 56192.55
$ws.Range("N131").Value = -66272.54999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4419.9
$ws.Range("I70").Value = 4599.8335
$ws.Range("J70").Value = 4150
$ws.Range("K70").Value = 4599.8335
$ws.Range("L70").Value = 4150
$ws.Range("M70").Value = -4329.8335
$ws.Range("N70").Value = -4690
$ws.Range("H73").Value = 4419.9
$ws.Range("I73").Value = 4599.8335
$ws.Range("J73").Value = 4150
$ws.Range("K73").Value = 4599.8335
$ws.Range("L73").Value = 4150
$ws.Range("M73").Value = -3663.8335
$ws.Range("N73").Value = -6022
$ws.Range("H102").Value = 2907.6667
$ws.Range("J102").Value = 1654
$ws.Range("L102").Value = 1654
$ws.Range("N102").Value = -4898
$ws.Range("H113").Value = 1939.8572
$ws.Range("I113").Value = 1926.3334
$ws.Range("J113").Value = 1950
$ws.Range("K113").Value = 1926.3334
$ws.Range("L113").Value = 1950
$ws.Range("M113").Value = 243.6666
$ws.Range("N113").Value = -6290
$ws.Range("H122").Value = 2994.75
$ws.Range("J122").Value = 3972
$ws.Range("L122").Value = 11916
$ws.Range("N122").Value = -16816
$ws.Range("H126").Value = 2528368.5
$ws.Range("I126").Value = 4632931
$ws.Range("J126").Value = 2893.7
$ws.Range("K126").Value = 13898793
$ws.Range("L126").Value = 8681.099999999999
$ws.Range("M126").Value = -13896323
$ws.Range("N126").Value = -13621.1
$ws.Range("H132").Value = 1042098.6
$ws.Range("I132").Value = 1540212.5
$ws.Range("K132").Value = 4620637.5
$ws.Range("M132").Value = -4618107.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3125.389
$ws.Range("I61").Value = 2189.9167
$ws.Range("K61").Value = 2189.9167
$ws.Range("M61").Value = -1987.9167
$ws.Range("H113").Value = 3125.389
$ws.Range("I113").Value = 2189.9167
$ws.Range("K113").Value = 2189.9167
$ws.Range("M113").Value = -19.91670000000022
$ws.Range("H136").Value = 4765.8945
$ws.Range("I136").Value = 2714.818
$ws.Range("J136").Value = 7586.125
$ws.Range("K136").Value = 8144.454000000001
$ws.Range("L136").Value = 22758.375
$ws.Range("M136").Value = -5594.454000000001
$ws.Range("N136").Value = -27858.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 446.85
$ws.Range("I113").Value = 293.9091
$ws.Range("J113").Value = 633.7778
$ws.Range("K113").Value = 881.7273
$ws.Range("L113").Value = 1901.3334
$ws.Range("M113").Value = 1288.2727
$ws.Range("N113").Value = -6241.3334
$ws.Range("H126").Value = 4472.8696
$ws.Range("I126").Value = 3681.625
$ws.Range("J126").Value = 6281.4287
$ws.Range("K126").Value = 11044.875
$ws.Range("L126").Value = 18844.2861
$ws.Range("M126").Value = -8574.875
$ws.Range("N126").Value = -23784.2861
$ws.Range("H136").Value = 22225050
$ws.Range("J136").Value = 3806.25
$ws.Range("L136").Value = 11418.75
$ws.Range("N136").Value = -16518.75

Write-Output "Applied Tonberry_Profits price updates"